$d = $word.ActiveDocument

# 1. "Building ... to be used for ..." -> "Built ... used for ..." (past tense).
$d.Content.Find.Execute( `
    "Building a neural network to be used for optical character recognition", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Built a neural network used for optical character recognition", 2)

# 2. Append the "with cross-entropy loss" detail right after "activation
#    functions", at the very end of that bullet's text. A throwaway marker
#    token is appended after it so the insertion point used for re-homing
#    the "_GoBack" bookmark below lands *inside* the paragraph's text run
#    rather than exactly on the paragraph mark, then the marker is removed
#    once the bookmark is anchored.
$d.Content.Find.Execute( `
    "activation functions", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "activation functions with cross-entropy lossZZMARKERZZ", 2)

$markerRange = $d.Content
$markerRange.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# 3. Word keeps a "_GoBack" bookmark at the location of the most recent
#    edit; re-home it here, immediately after the newly typed text.
$goBackRange = $d.Range($markerRange.Start, $markerRange.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# 4. Remove the throwaway marker now that the bookmark is anchored.
$markerRange2 = $d.Content
$markerRange2.Find.Execute("ZZMARKERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange2.Delete()
